# Insert a new sensor row (BMP180) right above the existing BMP280 row
# (old row 16 "BMP280" -> new row 17 "BMP280"; rows 16-20 shift to 17-21).
#
# Strategy: insert a blank row at 16 by shifting only columns A:L (so the
# sheet's used range/dimension stays A1:L21 instead of growing to the full
# row width), then copy the formatting from the row that is about to become
# row 17 (the old row 16) into the new row 16 so the new row picks up the
# same style indices (borders etc.) Excel would apply, and finally fill in
# the BMP180 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift A16:L20 down to A17:L21, leaving a blank A16:L16 behind.
$ws.Range("A16:L16").Insert(-4121)

# Copy formatting from the row right below (old row 16, now row 17) into
# the freshly inserted row 16, so the new row's cell styles match what
# Excel would naturally carry over for a row inserted in this table.
$ws.Range("A17:L17").Copy()
$ws.Range("A16:L16").PasteSpecial(-4122)

# Fill in the new BMP180 sensor's data.
$ws.Range("A16").Value = "BMP180"
$ws.Range("B16").Value = "1.8...3.6V"
$ws.Range("C16").Value = "I2C"
$ws.Range("D16").Value = "-40...85 °C"
$ws.Range("E16").Value = "±0.5 °C"
$ws.Range("F16").Value = "0.1 °C"
$ws.Range("J16").Value = "300...1100 hPA"
$ws.Range("K16").Value = "±1.0 hPa"
$ws.Range("L16").Value = "0.01 hPa"

# B16 ends up with no explicit cell style in the target workbook (unlike
# every other cell in column B, which carries style index 1).
$ws.Range("B16").ClearFormats()

# Match the author's final selection.
$ws.Range("L17").Select() | Out-Null
